$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: "All" -> "$ bold('All')"
$ws.Range("B1").Value = '$ bold(''All'')'

# C1: "European Union" -> "$ bold('European Union')"
$ws.Range("C1").Value = '$ bold(''European Union'')'

# A5: add atop()/bold() formatting with extra leading spaces and newline
$ws.Range("A5").Value = '$ atop(''                    Supports the GCS if its other members* cover 64-72% of world emissions'', 
                               ''*'' * bold(High) * '': Global South + China + EU + various HICs (UK, Japan, South Korea, Canada...)'')'

# A6: add atop()/bold() formatting with extra leading spaces and newline
$ws.Range("A6").Value = '$ atop(''     Supports the GCS if its other members* cover 64-72% of world emissions'',          
                                     ''*'' * bold(''High color'') * '': High + Distributive effects displayed using colors on world map'')'

# A7: add atop()/bold() formatting with extra leading spaces and newline
$ws.Range("A7").Value = '$ atop(''Supports the GCS if its other members* cover 56% of world emissions'', 
                              ''                                                                   *'' * bold(''Mid'') * '': Global South + China'')'

# A8: add atop()/bold() formatting with extra leading spaces and newline
$ws.Range("A8").Value = '$ atop(''Supports the GCS if its other members* cover 25-33% of world emissions'', 
                              ''                                                                            *'' * bold(''Low'') * '': Global South + EU'')'

